$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1,1).Value = "Datos actualizados a 9 de Octubre de 2020 a las 02:52"

# Row 4
$ws.Cells.Item(4,2).Value = 7831478
$ws.Cells.Item(4,3).Value = 54367
$ws.Cells.Item(4,4).Value = 5016139
$ws.Cells.Item(4,5).Value = 2597678
$ws.Cells.Item(4,7).Value = 880
$ws.Cells.Item(4,8).Value = 217661

# Row 26
$ws.Cells.Item(26,2).Value = 315514
$ws.Cells.Item(26,3).Value = 4401
$ws.Cells.Item(26,4).Value = 269500
$ws.Cells.Item(26,5).Value = 36347

# Row 29
$ws.Cells.Item(29,2).Value = 175559
$ws.Cells.Item(29,3).Value = 2436
$ws.Cells.Item(29,4).Value = 147508
$ws.Cells.Item(29,5).Value = 18494
$ws.Cells.Item(29,7).Value = 16
$ws.Cells.Item(29,8).Value = 9557

# Row 55
$ws.Cells.Item(55,1).Value = "Venezuela"
$ws.Cells.Item(55,2).Value = 81019
$ws.Cells.Item(55,3).Value = 615
$ws.Cells.Item(55,4).Value = 72196
$ws.Cells.Item(55,5).Value = 8145
$ws.Cells.Item(55,7).Value = 7
$ws.Cells.Item(55,8).Value = 678

# Row 56
$ws.Cells.Item(56,1).Value = "Honduras"
$ws.Cells.Item(56,2).Value = 81016
$ws.Cells.Item(56,3).Value = 354
$ws.Cells.Item(56,4).Value = 30590
$ws.Cells.Item(56,5).Value = 47960
$ws.Cells.Item(56,7).Value = 19
$ws.Cells.Item(56,8).Value = 2466

# Row 94
$ws.Cells.Item(94,2).Value = 15301
$ws.Cells.Item(94,3).Value = 77
$ws.Cells.Item(94,4).Value = 14365
$ws.Cells.Item(94,5).Value = 601

# Row 120
$ws.Cells.Item(120,1).Value = "Angola"
$ws.Cells.Item(120,2).Value = 5958
$ws.Cells.Item(120,3).Value = 95
$ws.Cells.Item(120,4).Value = 2635
$ws.Cells.Item(120,5).Value = 3115
$ws.Cells.Item(120,7).Value = 2
$ws.Cells.Item(120,8).Value = 208

# Row 121
$ws.Cells.Item(121,1).Value = "Cuba"
$ws.Cells.Item(121,2).Value = 5917
$ws.Cells.Item(121,3).Value = 19
$ws.Cells.Item(121,4).Value = 5371
$ws.Cells.Item(121,5).Value = 423
$ws.Cells.Item(121,8).Value = 123

# Row 122
$ws.Cells.Item(122,1).Value = "Malaui"
$ws.Cells.Item(122,2).Value = 5809
$ws.Cells.Item(122,3).Value = 6
$ws.Cells.Item(122,4).Value = 4626
$ws.Cells.Item(122,5).Value = 1003
$ws.Cells.Item(122,8).Value = 180

# Row 125
$ws.Cells.Item(125,4).Value = 5355
$ws.Cells.Item(125,5).Value = 7

# Row 129
$ws.Cells.Item(129,2).Value = 5062
$ws.Cells.Item(129,3).Value = 10
$ws.Cells.Item(129,5).Value = 85

# Row 130
$ws.Cells.Item(130,2).Value = 5004
$ws.Cells.Item(130,3).Value = 25
$ws.Cells.Item(130,4).Value = 4794
$ws.Cells.Item(130,5).Value = 104

# Row 133
$ws.Cells.Item(133,2).Value = 4853
$ws.Cells.Item(133,3).Value = 1
$ws.Cells.Item(133,5).Value = 2877

# Row 167
$ws.Cells.Item(167,2).Value = 1201
$ws.Cells.Item(167,3).Value = 1
$ws.Cells.Item(167,5).Value = 10

# Row 169
$ws.Cells.Item(169,2).Value = 921
$ws.Cells.Item(169,3).Value = 7
$ws.Cells.Item(169,4).Value = 890
$ws.Cells.Item(169,5).Value = 16

# Row 173
$ws.Cells.Item(173,4).Value = 671
$ws.Cells.Item(173,5).Value = 18

# Row 190
$ws.Cells.Item(190,2).Value = 214
$ws.Cells.Item(190,3).Value = 1
$ws.Cells.Item(190,4).Value = 211
